$wb = $excel.ActiveWorkbook

# The workbook's shared-string "Ready for handoff" is used as the localization
# status for both target languages, in three places:
#   - Overview sheet: zh-cn column (E) and de-de column (F), data rows 2-3
#   - zh-cn sheet: Status column (C), data rows 2-3
#   - de-de sheet: Status column (C), data rows 2-3
# All of these become "In Translation".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# The Status columns are narrower now that the text is shorter
# ("In Translation" vs "Ready for handoff"): Overview!E:F and the Status
# column (C) on the zh-cn / de-de sheets.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
